$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2 through 45: update date serial 45203 -> 45205
# (2023-10-04 -> 2023-10-06), keeping existing cell style/format intact.
for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45205
    }
}
